$d = $word.ActiveDocument

# Locate the target paragraph ("end: reset camera angle") by its text,
# then widen the range to the whole paragraph (including its end-of-
# paragraph mark) so InsertXML below replaces the paragraph in full
# rather than just the matched characters.
$found = $d.Content
$null = $found.Find.Execute("end: reset camera angle", $true, $false, $false,
                             $false, $false, $true, 1, $false, "", 0)
$target = $found.Paragraphs(1)
$r = $target.Range

$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$xml = "<w:p $ns>" +
         '<w:proofErr w:type="gramStart"/>' +
         '<w:r><w:t>end</w:t></w:r>' +
         '<w:proofErr w:type="gramEnd"/>' +
         '<w:r><w:t>: reset camera angle</w:t></w:r>' +
       '</w:p>' +
       "<w:p $ns/>" +
       "<w:p $ns>" +
         '<w:pPr><w:rPr><w:b/></w:rPr></w:pPr>' +
         '<w:r><w:rPr><w:b/></w:rPr><w:t>Toggle robot components:</w:t></w:r>' +
       '</w:p>' +
       "<w:p $ns>" +
         '<w:proofErr w:type="spellStart"/>' +
         '<w:proofErr w:type="gramStart"/>' +
         '<w:r><w:t>i</w:t></w:r>' +
         '<w:proofErr w:type="spellEnd"/>' +
         '<w:proofErr w:type="gramEnd"/>' +
         '<w:r><w:t>: cycle components</w:t></w:r>' +
       '</w:p>' +
       "<w:p $ns>" +
         '<w:proofErr w:type="gramStart"/>' +
         '<w:r><w:t>j</w:t></w:r>' +
         '<w:proofErr w:type="gramEnd"/>' +
         '<w:r><w:t>: turn component off</w:t></w:r>' +
       '</w:p>' +
       "<w:p $ns>" +
         '<w:r><w:t>u: turn component on</w:t></w:r>' +
       '</w:p>'

$null = $r.InsertXML($xml)
